$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 8
$ws.Range("G4").Value = 8
$ws.Range("G5").Value = 10
$ws.Range("G6").Value = 8
$ws.Range("G7").Value = 5
$ws.Range("G8").Value = 5
$ws.Range("G9").Value = 5
$ws.Range("G10").Value = 4
$ws.Range("G11").Value = 7
$ws.Range("G12").Value = 10
$ws.Range("G13").Value = 10
$ws.Range("G14").Value = 5
$ws.Range("G15").Value = 5
$ws.Range("G16").Value = 6
$ws.Range("G17").Value = 4
$ws.Range("G18").Value = 11
$ws.Range("G19").Value = 8
$ws.Range("G20").Value = 7
$ws.Range("G21").Value = 7
$ws.Range("G22").Value = 6
$ws.Range("G23").Value = 9
$ws.Range("G24").Value = 4
$ws.Range("G25").Value = 7
$ws.Range("G26").Value = 8
$ws.Range("G27").Value = 3
$ws.Range("G28").Value = 12
$ws.Range("G29").Value = 3
$ws.Range("G30").Value = 8
$ws.Range("G31").Value = 7
$ws.Range("G32").Value = 4
$ws.Range("G33").Value = 7
$ws.Range("G34").Value = 11
$ws.Range("G35").Value = 5
$ws.Range("G36").Value = 6
$ws.Range("G37").Value = 5
$ws.Range("G38").Value = 10
$ws.Range("G39").Value = 5
$ws.Range("G40").Value = 2
$ws.Range("G41").Value = 0
$ws.Range("G42").Value = 2
